# Apply the "first model for the path following controller" edit:
#   - keep Tabelle1 as-is, but update its view (zoom 130%, selection A25)
#   - append a new worksheet "Sheet1" after Tabelle1, becoming the active/tabbed sheet
#   - populate the new sheet with a small "ToDos" / "Aufgabe" table
#   - the B2/H2 header cells use a bold, 14pt font

$wb = $excel.ActiveWorkbook

# --- Tabelle1: update the view before we move focus away from it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$excel.ActiveWindow.Zoom = 130
$ws1.Range("A25").Select() | Out-Null

# --- add the new sheet right after Tabelle1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet1"

# --- cell values, entered column-by-column (B then H) to match the
#     original authoring order of the shared-string table ---
$ws2.Range("B2").Value = "ToDos"
$ws2.Range("B4").Value = "Code aus Aufgabe in Bericht oder in Ordner lassen?"
$ws2.Range("H2").Value = "Aufgabe"
$ws2.Range("H4").Value = "8.3.2"

# --- header row formatting (bold, 14pt) ---
$b2 = $ws2.Range("B2")
$b2.Font.Bold = $true
$b2.Font.Size = 14

# copy B2's formatting onto H2 so both share the same cell style (instead of
# creating a second, separate style entry)
$b2.Copy() | Out-Null
$ws2.Range("H2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws2.Range("I9").Select() | Out-Null
